$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply updated symbol list values cell by cell, row by row (per commit:
# "Updated symbol list on Thu Jan 12 14:34:18 UTC 2023 with GitHub Actions").
# Columns B (Coin) and C (Link) are plain text and are assigned directly.
# Columns D (Price) and E (Volume 1h) look numeric/percentage, so we force
# them to stay Text (NumberFormat "@") before assignment, matching the
# original workbook where every data cell is stored as a text string,
# then reset the style back to Normal so no stray formatting is left behind.

# Row 2
$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("D2").Value = '283.21'
$ws.Range("E2").Value = '2.09%'
$ws.Range("D2:E2").Style = "Normal"

# Row 3
$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("D3").Value = '28.44'
$ws.Range("E3").Value = '3.45%'
$ws.Range("D3:E3").Style = "Normal"

# Row 4
$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("D4").Value = '5.051'
$ws.Range("E4").Value = '4.23%'
$ws.Range("D4:E4").Style = "Normal"

# Row 5
$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("D5").Value = '0.06515'
$ws.Range("E5").Value = '2.27%'
$ws.Range("D5:E5").Style = "Normal"

# Row 6
$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("D6").Value = '7.249'
$ws.Range("E6").Value = '3.56%'
$ws.Range("D6:E6").Style = "Normal"

# Row 7
$ws.Range("B7").Value = 'GateToken'
$ws.Range("C7").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("D7:E7").NumberFormat = "@"
$ws.Range("D7").Value = '3.385'
$ws.Range("E7").Value = '2.56%'
$ws.Range("D7:E7").Style = "Normal"

# Row 8
$ws.Range("B8").Value = 'FTXToken'
$ws.Range("C8").Value = 'https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt'
$ws.Range("D8:E8").NumberFormat = "@"
$ws.Range("D8").Value = '1.412'
$ws.Range("E8").Value = '1.12%'
$ws.Range("D8:E8").Style = "Normal"

# Row 9
$ws.Range("B9").Value = 'MXToken'
$ws.Range("C9").Value = 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx'
$ws.Range("D9:E9").NumberFormat = "@"
$ws.Range("D9").Value = '0.9225'
$ws.Range("E9").Value = '5.36%'
$ws.Range("D9:E9").Style = "Normal"

# Row 10
$ws.Range("B10").Value = 'WazirX'
$ws.Range("C10").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("D10").Value = '0.1546'
$ws.Range("E10").Value = '1.99%'
$ws.Range("D10:E10").Style = "Normal"

# Row 11
$ws.Range("B11").Value = 'LiechtensteinCryptoassetsExchange'
$ws.Range("C11").Value = 'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx'
$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("D11").Value = '0.06451'
$ws.Range("E11").Value = '24.61%'
$ws.Range("D11:E11").Style = "Normal"

# Row 12
$ws.Range("B12").Value = 'MandalaExchangeToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx'
$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("D12").Value = '0.07617'
$ws.Range("E12").Value = '1.26%'
$ws.Range("D12:E12").Style = "Normal"

# Row 13
$ws.Range("B13").Value = 'BitrueCoin'
$ws.Range("C13").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("D13").Value = '0.02894'
$ws.Range("E13").Value = '-1.34%'
$ws.Range("D13:E13").Style = "Normal"

# Row 14
$ws.Range("B14").Value = 'BitMartToken'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("D14").Value = '0.08943'
$ws.Range("E14").Value = '-0.22%'
$ws.Range("D14:E14").Style = "Normal"

# Row 15
$ws.Range("B15").Value = 'BitForexToken'
$ws.Range("C15").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("D15").Value = '0.001586'
$ws.Range("E15").Value = '1.24%'
$ws.Range("D15:E15").Style = "Normal"

# Row 16
$ws.Range("B16").Value = 'One'
$ws.Range("C16").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0006404'
$ws.Range("E16").Value = '0.24%'
$ws.Range("D16:E16").Style = "Normal"

# Row 17
$ws.Range("B17").Value = 'TigerCash'
$ws.Range("C17").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '0.006030'
$ws.Range("D17").Style = "Normal"

# Row 18
$ws.Range("B18").Value = 'LEO'
$ws.Range("C18").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("D18:E18").NumberFormat = "@"
$ws.Range("D18").Value = '3.448'
$ws.Range("E18").Value = '-0.82%'
$ws.Range("D18:E18").Style = "Normal"

# Row 19
$ws.Range("D19:E19").NumberFormat = "@"
$ws.Range("D19").Value = '2.242'
$ws.Range("E19").Value = '-0.26%'
$ws.Range("D19:E19").Style = "Normal"

# Row 20
$ws.Range("D20:E20").NumberFormat = "@"
$ws.Range("D20").Value = '0.3192'
$ws.Range("E20").Value = '1.42%'
$ws.Range("D20:E20").Style = "Normal"

# Row 21
$ws.Range("D21:E21").NumberFormat = "@"
$ws.Range("D21").Value = '0.1289'
$ws.Range("E21").Value = '-2.64%'
$ws.Range("D21:E21").Style = "Normal"

# Row 22
$ws.Range("D22:E22").NumberFormat = "@"
$ws.Range("D22").Value = '3.991'
$ws.Range("E22").Value = '2.12%'
$ws.Range("D22:E22").Style = "Normal"

# Row 23
$ws.Range("D23:E23").NumberFormat = "@"
$ws.Range("D23").Value = '0.1541'
$ws.Range("E23").Value = '1.44%'
$ws.Range("D23:E23").Style = "Normal"

# Row 24
$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("D24").Value = '0.04416'
$ws.Range("E24").Value = '0.21%'
$ws.Range("D24:E24").Style = "Normal"

# Row 25
$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("D25").Value = '0.001186'
$ws.Range("E25").Value = '0.88%'
$ws.Range("D25:E25").Style = "Normal"

# Row 26
$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("D26").Value = '0.004410'
$ws.Range("E26").Value = '13.16%'
$ws.Range("D26:E26").Style = "Normal"

# Row 27
$ws.Range("D27:E27").NumberFormat = "@"
$ws.Range("D27").Value = '0.0001252'
$ws.Range("E27").Value = '6.02%'
$ws.Range("D27:E27").Style = "Normal"

# Row 28
$ws.Range("D28:E28").NumberFormat = "@"
$ws.Range("D28").Value = '0.0001618'
$ws.Range("E28").Value = '-1.60%'
$ws.Range("D28:E28").Style = "Normal"

# Row 40
$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("D40").Value = '0.04153'
$ws.Range("E40").Value = '2.08%'
$ws.Range("D40:E40").Style = "Normal"

# Row 41
$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("D41").Value = '0.006679'
$ws.Range("E41").Value = '-2.07%'
$ws.Range("D41:E41").Style = "Normal"

# Row 42
$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1227'
$ws.Range("E42").Value = '-13.20%'
$ws.Range("D42:E42").Style = "Normal"

# Row 43
$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("D43").Value = '0.002193'
$ws.Range("E43").Value = '16.05%'
$ws.Range("D43:E43").Style = "Normal"

# Row 44
$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("D44").Value = '0.01215'
$ws.Range("E44").Value = '4.10%'
$ws.Range("D44:E44").Style = "Normal"

# Row 45
$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("D45").Value = '0.00005632'
$ws.Range("E45").Value = '5.18%'
$ws.Range("D45:E45").Style = "Normal"

# Row 47
$ws.Range("D47:E47").NumberFormat = "@"
$ws.Range("D47").Value = '0.01302'
$ws.Range("E47").Value = '-29.64%'
$ws.Range("D47:E47").Style = "Normal"
